$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New row (17) goes right after the last existing data row (16).
# Borrow the border/fill formatting of an existing fully-bordered data
# row (row 11) so the new row's cells share the same border style used
# throughout the table instead of creating a brand new border definition.
$ws.Range("A11:E11").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column B ("Jira id") on this new row wraps its text, matching the
# other wrap-text styled cells in the sheet.
$ws.Range("B17").WrapText = $true

# Fill in the new test case data.
$ws.Range("A17").Value = "RCC100"
$ws.Range("B17").Value = "OBT"
$ws.Range("C17").Value = "Verify the sorting options"
$ws.Range("D17").Value = "Y"
$ws.Range("E17").Value = ""

# Match the selection Excel leaves behind after entering a new row of data.
$ws.Range("A17:E17").Select()
